$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (copy formatting from the existing header style in F1)
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Updated metric values (B, C, D) for rows 2-4
$ws.Range("B2").Value = 0.04896347446528382
$ws.Range("C2").Value = 0.9985600677500517
$ws.Range("D2").Value = 0.161481537436065

$ws.Range("B3").Value = 0.07415511211264601
$ws.Range("C3").Value = 0.9992966055291462
$ws.Range("D3").Value = 0.2077363080320377

$ws.Range("B4").Value = 0.09185755673884506
$ws.Range("C4").Value = 0.998756305383627
$ws.Range("D4").Value = 0.2381653281603272

# New Elapsed Time / CPU columns for rows 2-4
$ws.Range("G2").Value = 0.3776785511166963
$ws.Range("H2").Value = 0.968

$ws.Range("G3").Value = 0.3776785511166963
$ws.Range("H3").Value = 0.968

$ws.Range("G4").Value = 0.3776785511166963
$ws.Range("H4").Value = 0.968
